$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(28).Insert()

$ws.Range("A28").Value = 10
$ws.Range("B28").Value = "Vega Modelo de Temuco"
$ws.Range("C28").Value = "La Araucanía"
$ws.Range("D28").Value = 44414
$ws.Range("E28").Value = 9
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100104
$ws.Range("H28").Value = "Frutos de pepita"
$ws.Range("I28").Value = 100104003
$ws.Range("J28").Value = "Membrillo"
$ws.Range("K28").Value = "Champion"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 95
$ws.Range("N28").Value = 12000
$ws.Range("O28").Value = 13000
$ws.Range("P28").Value = 12526
$ws.Range("Q28").Value = "$/bandeja 18 kilos granel"
$ws.Range("R28").Value = "Región de O'Higgins"
$ws.Range("S28").Value = 696
$ws.Range("T28").Value = 18
